$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the spend-list text in B7 (small whitespace tightening between
#    "Mortgage" and "Utilities and telecommunications": 79 spaces -> 77 spaces).
$ws.Range("B7").Value = "•Food`n•Clothing`n•Household supplies and personal care                                 •Household items (TV, electronics, furniture, appliances)  `n•Recreational goods (sports and fitness equipment, bicycles, toys, games)                             `n•Rent`n•Mortgage                                                                             •Utilities and telecommunications`n•Vehicle payments                                                              •Paying down credit card, student loans, or other debts  `n•Charitable donations or giving to family members                             `n•Savings or other investments`n•Other, please describe:"

# 2. Apply top vertical alignment across the used grid (keeps existing
#    wrap/fonts, just adds vertical="top" to every in-use cell).
$ws.Range("A1:E10").VerticalAlignment = -4160
$ws.Range("A11:A16").VerticalAlignment = -4160

# 3. Row 9 height shrinks slightly (106 -> 105).
$ws.Rows.Item(9).RowHeight = 105

# 4. Final recorded selection is B10.
$ws.Range("B10").Select()
